$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 27902.732
$ws.Range("I86").Value = 1298
$ws.Range("J86").Value = 51181.875
$ws.Range("K86").Value = 1298
$ws.Range("L86").Value = 51181.875
$ws.Range("M86").Value = -175
$ws.Range("N86").Value = -53427.875
$ws.Range("H89").Value = 27902.732
$ws.Range("I89").Value = 1298
$ws.Range("J89").Value = 51181.875
$ws.Range("K89").Value = 6490
$ws.Range("L89").Value = 255909.375
$ws.Range("M89").Value = -874
$ws.Range("N89").Value = -267141.375
$ws.Range("H106").Value = 2092.2307
$ws.Range("I106").Value = 1521
$ws.Range("J106").Value = 3996.3333
$ws.Range("K106").Value = 1521
$ws.Range("L106").Value = 3996.3333
$ws.Range("M106").Value = -890
$ws.Range("N106").Value = -5258.3333
$ws.Range("H112").Value = 1127.8975
$ws.Range("J112").Value = 1145.091
$ws.Range("L112").Value = 3435.273
$ws.Range("N112").Value = -5651.272999999999
$ws.Range("H129").Value = 961.7027
$ws.Range("I129").Value = 391.92856
$ws.Range("J129").Value = 1308.5217
$ws.Range("K129").Value = 1175.78568
$ws.Range("L129").Value = 3925.5651
$ws.Range("M129").Value = 3824.21432
$ws.Range("N129").Value = -13925.5651
$ws.Range("H137").Value = 1339.973
$ws.Range("I137").Value = 1396
$ws.Range("J137").Value = 1280.8334
$ws.Range("K137").Value = 4188
$ws.Range("L137").Value = 3842.5002
$ws.Range("M137").Value = -1638
$ws.Range("N137").Value = -8942.5002
$ws.Range("H138").Value = 16669926
$ws.Range("I138").Value = 1356.7
$ws.Range("J138").Value = 50007064
$ws.Range("K138").Value = 4070.1
$ws.Range("L138").Value = 150021192
$ws.Range("M138").Value = 1069.9
$ws.Range("N138").Value = -150031472

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 27348.545
$ws.Range("I5").Value = 22303.555
$ws.Range("K5").Value = 22303.555
$ws.Range("M5").Value = -22191.555
$ws.Range("H32").Value = 6133.913
$ws.Range("I32").Value = 4054
$ws.Range("K32").Value = 4054
$ws.Range("M32").Value = -3767
$ws.Range("H74").Value = 1446.9688
$ws.Range("I74").Value = 1429.8077
$ws.Range("K74").Value = 1429.8077
$ws.Range("M74").Value = -555.8077000000001
$ws.Range("H77").Value = 1446.9688
$ws.Range("I77").Value = 1429.8077
$ws.Range("K77").Value = 7149.038500000001
$ws.Range("M77").Value = -2781.038500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 27348.545
$ws.Range("I4").Value = 22303.555
$ws.Range("K4").Value = 22303.555
$ws.Range("M4").Value = -22188.555
$ws.Range("H105").Value = 1800
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5094

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 195.63637
$ws.Range("J7").Value = 195.63637
$ws.Range("L7").Value = 195.63637
$ws.Range("N7").Value = -421.63637
$ws.Range("H99").Value = 3211803.5
$ws.Range("I99").Value = 4573719.5
$ws.Range("J99").Value = 33999.668
$ws.Range("K99").Value = 4573719.5
$ws.Range("L99").Value = 33999.668
$ws.Range("M99").Value = -4572221.5
$ws.Range("N99").Value = -36995.668
$ws.Range("H126").Value = 3211803.5
$ws.Range("I126").Value = 4573719.5
$ws.Range("J126").Value = 33999.668
$ws.Range("K126").Value = 13721158.5
$ws.Range("L126").Value = 101999.004
$ws.Range("M126").Value = -13718688.5
$ws.Range("N126").Value = -106939.004
$ws.Range("H132").Value = 1975.6842
$ws.Range("I132").Value = 1610.129
$ws.Range("J132").Value = 3594.5715
$ws.Range("K132").Value = 4830.387
$ws.Range("L132").Value = 10783.7145
$ws.Range("M132").Value = -2300.387
$ws.Range("N132").Value = -15843.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 137.53572
$ws.Range("J33").Value = 241.33333
$ws.Range("L33").Value = 1447.99998
$ws.Range("N33").Value = -2013.99998
$ws.Range("H35").Value = 1325
$ws.Range("H68").Value = 974.1340300000001
$ws.Range("I68").Value = 778.3095
$ws.Range("J68").Value = 1123.6727
$ws.Range("K68").Value = 2334.9285
$ws.Range("L68").Value = 3371.0181
$ws.Range("M68").Value = -1523.9285
$ws.Range("N68").Value = -4993.0181
$ws.Range("H71").Value = 974.1340300000001
$ws.Range("I71").Value = 778.3095
$ws.Range("J71").Value = 1123.6727
$ws.Range("K71").Value = 7004.7855
$ws.Range("L71").Value = 10113.0543
$ws.Range("M71").Value = -2948.7855
$ws.Range("N71").Value = -18225.0543
$ws.Range("H122").Value = 511
$ws.Range("I122").Value = 360.34784
$ws.Range("J122").Value = 1666
$ws.Range("K122").Value = 3243.13056
$ws.Range("L122").Value = 14994
$ws.Range("M122").Value = -793.1305600000001
$ws.Range("N122").Value = -19894
$ws.Range("H138").Value = 1494
$ws.Range("I138").Value = 988
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 2964
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 2176
$ws.Range("N138").Value = -16280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 146
$ws.Range("I2").Value = 42
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 42
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 71
$ws.Range("N2").Value = -476
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66095.625
$ws.Range("I7").Value = 86960.836
$ws.Range("K7").Value = 86960.836
$ws.Range("M7").Value = -86848.836
$ws.Range("H22").Value = 1199.4783
$ws.Range("I22").Value = 691.53845
$ws.Range("J22").Value = 1859.8
$ws.Range("K22").Value = 691.53845
$ws.Range("L22").Value = 1859.8
$ws.Range("M22").Value = -396.53845
$ws.Range("N22").Value = -2449.8
$ws.Range("H27").Value = 1199.4783
$ws.Range("I27").Value = 691.53845
$ws.Range("J27").Value = 1859.8
$ws.Range("K27").Value = 691.53845
$ws.Range("L27").Value = 1859.8
$ws.Range("M27").Value = -584.53845
$ws.Range("N27").Value = -2073.8
$ws.Range("H40").Value = 60621.777
$ws.Range("I40").Value = 252000
$ws.Range("J40").Value = 5942.2856
$ws.Range("K40").Value = 252000
$ws.Range("L40").Value = 5942.2856
$ws.Range("M40").Value = -251864
$ws.Range("N40").Value = -6214.2856
$ws.Range("H87").Value = 10189
$ws.Range("J87").Value = 10189
$ws.Range("L87").Value = 10189
$ws.Range("N87").Value = -12435
$ws.Range("H90").Value = 10189
$ws.Range("J90").Value = 10189
$ws.Range("L90").Value = 30567
$ws.Range("N90").Value = -41799
$ws.Range("H122").Value = 27780030
$ws.Range("I122").Value = 37039036
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 111117108
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -111114658
$ws.Range("N122").Value = -13915
$ws.Range("H126").Value = 66095.625
$ws.Range("I126").Value = 86960.836
$ws.Range("K126").Value = 260882.508
$ws.Range("M126").Value = -258412.508
$ws.Range("H132").Value = 4350.1396
$ws.Range("I132").Value = 4865.476
$ws.Range("J132").Value = 3858.2273
$ws.Range("K132").Value = 14596.428
$ws.Range("L132").Value = 11574.6819
$ws.Range("M132").Value = -12066.428
$ws.Range("N132").Value = -16634.6819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2416.7874
$ws.Range("I132").Value = 1709.7037
$ws.Range("J132").Value = 3371.35
$ws.Range("K132").Value = 5129.1111
$ws.Range("L132").Value = 10114.05
$ws.Range("M132").Value = -2599.1111
$ws.Range("N132").Value = -15174.05
